$wb = $excel.ActiveWorkbook

# --- 1) "Actions" sheet: row 4, column D - wizard name changes from
#        "On Instance" wording to the existing "On Entity Type" wording.
$wsActions = $wb.Worksheets.Item("Actions")
$wsActions.Cells.Item(4, 4).Value = "[Wizard Editor] Wizard Execution On Entity Type"
$wsActions.Range("D4").Select()

# --- 2) "Activities in Page" sheet: add a new row for the
#        "PuntoDiPrelievo" page tied to the ScheduleWorkflowForecast activity.
$wsActInPage = $wb.Worksheets.Item("Activities in Page")
$wsActInPage.Cells.Item(4, 1).Value = "CREATE/MODIFY"
$wsActInPage.Cells.Item(4, 2).Value = "ScheduleWorkflowForecast"
$wsActInPage.Cells.Item(4, 3).Value = "PuntoDiPrelievo"
$wsActInPage.Cells.Item(4, 4).Value = 100
$wsActInPage.Range("A4:C4").Select()

# --- 3) "Permissions" sheet: add a new row granting the intesa_user
#        group access to the new PuntoDiPrelievo page.
$wsPermissions = $wb.Worksheets.Item("Permissions")
$wsPermissions.Cells.Item(5, 1).Value = "CREATE/MODIFY"
$wsPermissions.Cells.Item(5, 2).Value = "ScheduleWorkflowForecast"
$wsPermissions.Cells.Item(5, 3).Value = "PuntoDiPrelievo"
$wsPermissions.Cells.Item(5, 4).Value = "intesa_user"
$wsPermissions.Range("D4:D5").Select()
